$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels: BP1/BQ1 (average_doctor <-> average_doctor_old)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4
$ws.Range("E4").Value = 0.424
$ws.Range("F4").Value = 0.07099999999999999
$ws.Range("G4").Value = 0.266
$ws.Range("N4").Value = 0.422
$ws.Range("O4").Value = 0.06
$ws.Range("P4").Value = 0.245
$ws.Range("Q4").Value = 0.025
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.132
$ws.Range("W4").Value = 0.295
$ws.Range("X4").Value = 0.111
$ws.Range("Y4").Value = 0.333
$ws.Range("AI4").Value = 0.297
$ws.Range("AJ4").Value = 0.08799999999999999
$ws.Range("AK4").Value = 0.296
$ws.Range("AU4").Value = 0.191
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.17
$ws.Range("BA4").Value = 1.982
$ws.Range("BB4").Value = 0.154
$ws.Range("BC4").Value = 0.393
$ws.Range("BG4").Value = 0.723
$ws.Range("BH4").Value = 0.145
$ws.Range("BI4").Value = 0.381
$ws.Range("BM4").Value = 0.711
$ws.Range("BN4").Value = 0.079
$ws.Range("BO4").Value = 0.281
$ws.Range("BP4").Value = 0.661
$ws.Range("BQ4").Value = 0.704

# Row 5
$ws.Range("E5").Value = 0.538
$ws.Range("F5").Value = 0.08699999999999999
$ws.Range("G5").Value = 0.295
$ws.Range("N5").Value = 0.743
$ws.Range("O5").Value = 0.08
$ws.Range("P5").Value = 0.284
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.08400000000000001
$ws.Range("W5").Value = 0.283
$ws.Range("X5").Value = 0.111
$ws.Range("Y5").Value = 0.333
$ws.Range("AI5").Value = 0.315
$ws.Range("AJ5").Value = 0.094
$ws.Range("AK5").Value = 0.307
$ws.Range("AU5").Value = 0.366
$ws.Range("AV5").Value = 0.091
$ws.Range("AW5").Value = 0.301
$ws.Range("BA5").Value = 1.331
$ws.Range("BB5").Value = 0.081
$ws.Range("BC5").Value = 0.284
$ws.Range("BG5").Value = 0.389
$ws.Range("BH5").Value = 0.049
$ws.Range("BI5").Value = 0.221
$ws.Range("BM5").Value = 0.549
$ws.Range("BN5").Value = 0.064
$ws.Range("BO5").Value = 0.252
$ws.Range("BP5").Value = 0.444
$ws.Range("BQ5").Value = 0.457

# Row 6
$ws.Range("E6").Value = 0.474
$ws.Range("N6").Value = 0.538
$ws.Range("Q6").Value = 0.02
$ws.Range("W6").Value = 0.289
$ws.Range("AI6").Value = 0.306
$ws.Range("AU6").Value = 0.251
$ws.Range("BA6").Value = 1.584
$ws.Range("BG6").Value = 0.506
$ws.Range("BM6").Value = 0.62
$ws.Range("BP6").Value = 0.528
$ws.Range("BQ6").Value = 0.551

# Row 7
$ws.Range("E7").Value = 0.511
$ws.Range("N7").Value = 0.645
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.285
$ws.Range("AI7").Value = 0.311
$ws.Range("AU7").Value = 0.309
$ws.Range("BA7").Value = 1.421
$ws.Range("BG7").Value = 0.429
$ws.Range("BM7").Value = 0.575
$ws.Range("BP7").Value = 0.474
$ws.Range("BQ7").Value = 0.49

# Row 8
$ws.Range("E8").Value = 0.594
$ws.Range("F8").Value = 0.115
$ws.Range("G8").Value = 0.34
$ws.Range("N8").Value = 0.771
$ws.Range("O8").Value = 0.06900000000000001
$ws.Range("P8").Value = 0.263
$ws.Range("Q8").Value = 0.018
$ws.Range("S8").Value = 0.11
$ws.Range("W8").Value = 0.308
$ws.Range("AI8").Value = 0.338
$ws.Range("AJ8").Value = 0.127
$ws.Range("AK8").Value = 0.356
$ws.Range("AU8").Value = 0.309
$ws.Range("AV8").Value = 0.08400000000000001
$ws.Range("AW8").Value = 0.29
$ws.Range("BA8").Value = 1.733
$ws.Range("BB8").Value = 0.125
$ws.Range("BC8").Value = 0.353
$ws.Range("BG8").Value = 0.555
$ws.Range("BH8").Value = 0.109
$ws.Range("BI8").Value = 0.33
$ws.Range("BM8").Value = 0.6889999999999999
$ws.Range("BN8").Value = 0.06900000000000001
$ws.Range("BO8").Value = 0.262
$ws.Range("BP8").Value = 0.578
$ws.Range("BQ8").Value = 0.601

# Row 9
$ws.Range("E9").Value = 0.528
$ws.Range("F9").Value = 0.249
$ws.Range("G9").Value = 0.499
$ws.Range("N9").Value = 0.663
$ws.Range("O9").Value = 0.223
$ws.Range("P9").Value = 0.473
$ws.Range("W9").Value = 0.202
$ws.Range("X9").Value = 0.161
$ws.Range("Y9").Value = 0.402
$ws.Range("AI9").Value = 0.258
$ws.Range("AJ9").Value = 0.192
$ws.Range("AK9").Value = 0.438
$ws.Range("BA9").Value = 1.673
$ws.Range("BB9").Value = 0.247
$ws.Range("BC9").Value = 0.497
$ws.Range("BG9").Value = 0.584
$ws.Range("BH9").Value = 0.243
$ws.Range("BI9").Value = 0.493
$ws.Range("BM9").Value = 0.64
$ws.Range("BN9").Value = 0.23
$ws.Range("BO9").Value = 0.48
$ws.Range("BP9").Value = 0.5580000000000001
$ws.Range("BQ9").Value = 0.579

# Row 10
$ws.Range("E10").Value = 0.663
$ws.Range("F10").Value = 0.223
$ws.Range("G10").Value = 0.473
$ws.Range("N10").Value = 0.865
$ws.Range("O10").Value = 0.117
$ws.Range("P10").Value = 0.342
$ws.Range("W10").Value = 0.382
$ws.Range("X10").Value = 0.236
$ws.Range("Y10").Value = 0.486
$ws.Range("AI10").Value = 0.371
$ws.Range("AJ10").Value = 0.233
$ws.Range("AK10").Value = 0.483
$ws.Range("AU10").Value = 0.303
$ws.Range("AV10").Value = 0.211
$ws.Range("AW10").Value = 0.46
$ws.Range("BA10").Value = 2.055
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.64
$ws.Range("BH10").Value = 0.23
$ws.Range("BI10").Value = 0.48
$ws.Range("BM10").Value = 0.831
$ws.Range("BN10").Value = 0.14
$ws.Range("BO10").Value = 0.374
$ws.Range("BP10").Value = 0.6850000000000001
$ws.Range("BQ10").Value = 0.719

# Row 11
$ws.Range("E11").Value = 0.697
$ws.Range("F11").Value = 0.211
$ws.Range("G11").Value = 0.46
$ws.Range("N11").Value = 0.888
$ws.Range("O11").Value = 0.1
$ws.Range("P11").Value = 0.316
$ws.Range("W11").Value = 0.382
$ws.Range("X11").Value = 0.236
$ws.Range("Y11").Value = 0.486
$ws.Range("AI11").Value = 0.404
$ws.Range("AJ11").Value = 0.241
$ws.Range("AK11").Value = 0.491
$ws.Range("AU11").Value = 0.427
$ws.Range("AV11").Value = 0.245
$ws.Range("AW11").Value = 0.495
$ws.Range("BA11").Value = 2.055
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.64
$ws.Range("BH11").Value = 0.23
$ws.Range("BI11").Value = 0.48
$ws.Range("BM11").Value = 0.831
$ws.Range("BN11").Value = 0.14
$ws.Range("BO11").Value = 0.374
$ws.Range("BP11").Value = 0.6850000000000001
$ws.Range("BQ11").Value = 0.722

# Row 12
$ws.Range("E12").Value = 1.435
$ws.Range("F12").Value = 0.794
$ws.Range("G12").Value = 0.891
$ws.Range("N12").Value = 1.494
$ws.Range("O12").Value = 1.089
$ws.Range("P12").Value = 1.044
$ws.Range("W12").Value = 1.647
$ws.Range("X12").Value = 0.581
$ws.Range("Y12").Value = 0.762
$ws.Range("AI12").Value = 1.722
$ws.Range("AJ12").Value = 1.312
$ws.Range("AK12").Value = 1.145
$ws.Range("AU12").Value = 2.725
$ws.Range("AV12").Value = 2.699
$ws.Range("AW12").Value = 1.643
$ws.Range("BA12").Value = 3.716
$ws.Range("BB12").Value = 0.412
$ws.Range("BC12").Value = 0.642
$ws.Range("BG12").Value = 1.105
$ws.Range("BH12").Value = 0.129
$ws.Range("BI12").Value = 0.36
$ws.Range("BM12").Value = 1.284
$ws.Range("BN12").Value = 0.311
$ws.Range("BO12").Value = 0.5580000000000001
$ws.Range("BP12").Value = 1.239
$ws.Range("BQ12").Value = 1.262

# Row 13
$ws.Range("E13").Value = 1.546
$ws.Range("F13").Value = 0.526
$ws.Range("G13").Value = 0.725
$ws.Range("N13").Value = 2.103
$ws.Range("O13").Value = 0.901
$ws.Range("P13").Value = 0.949
$ws.Range("W13").Value = 1.022
$ws.Range("X13").Value = 0.187
$ws.Range("Y13").Value = 0.432
$ws.Range("AI13").Value = 1.265
$ws.Range("AJ13").Value = 0.376
$ws.Range("AK13").Value = 0.613
$ws.Range("AU13").Value = 2.272
$ws.Range("AV13").Value = 0.921
$ws.Range("AW13").Value = 0.959
$ws.Range("BA13").Value = 2.349
$ws.Range("BB13").Value = 0.303
$ws.Range("BC13").Value = 0.55
$ws.Range("BG13").Value = 0.575
$ws.Range("BH13").Value = 0.05
$ws.Range("BI13").Value = 0.225
$ws.Range("BM13").Value = 0.888
$ws.Range("BN13").Value = 0.238
$ws.Range("BO13").Value = 0.488
$ws.Range("BP13").Value = 0.783
$ws.Range("BQ13").Value = 0.732
